# Reorder the comma-separated "Recorded By" names in column G so that
# "System" always appears first in the list, matching the canonical
# edit recorded for this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -notmatch ",") {
        continue
    }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) {
        continue
    }

    if ($parts[0] -ne "System") {
        $reversed = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
